$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Update the "GB" (data type indicator) column for the 4 display screens'
# numeric-indicator rows: change from "N" to "<value>"
$ws.Range("E28").Value = "<value>"
$ws.Range("E63").Value = "<value>"
$ws.Range("E98").Value = "<value>"
$ws.Range("E133").Value = "<value>"

# Add 4 new "display gare" rows (150-153), following the same pattern as
# the other "Extra" rows in the sheet (e.g. row 149's neighbours use
# TypographyName "Extra", Alignment "Left", GB "N", Direction "LTR")
$rows = @(150, 151, 152, 153)
$ids = @("SingleUseId148", "SingleUseId149", "SingleUseId150", "SingleUseId151")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Range("B$r").Value = $ids[$i]
    $ws.Range("C$r").Value = "Extra"
    $ws.Range("D$r").Value = "Left"
    $ws.Range("E$r").Value = "N"
    $ws.Range("F$r").Value = "LTR"
    # Newly populated cells otherwise inherit the column's default style;
    # reset to Normal so they match plain, unstyled cells like their peers.
    $ws.Range("B$r`:F$r").Style = "Normal"
}
